$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new weekly progress row (row 35) - week of 6/12/2025
$ws.Range("D35").Value = Get-Date -Year 2025 -Month 12 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Range("E35").Value = 0
$ws.Range("F35").Value = 518
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 1012
$ws.Range("J35").Value = "N/A"

# Move the active selection as it was left after editing
$ws.Range("I37").Select()

$wb.Save()
